# "_06_WriteInTheExcel" - keep appending the same text to the next empty
# row in column A (instead of always overwriting A1).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$text = "Merhaba Dünya"

# Find the first empty row below the current used range in column A and
# append three more rows with the same value, same as calling the
# "write to next row" helper three times.
for ($i = 1; $i -le 3; $i++) {
    $nextRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row + 1
    $ws.Cells.Item($nextRow, 1).Value = $text
}
